$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking strings (e.g. "320.27") remain plain text,
# matching the inline-string storage used in the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.464.28"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "1.845.06"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  +2.53%  "
$ws.Range("D5").Value = "320.27"
$ws.Range("E5").Value = "  +2.58%  "
$ws.Range("D6").Value = "1.028"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("D7").Value = "0.4354"
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("D8").Value = "0.3751"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").Value = "0.07355"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").Value = "0.8696"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").Value = "1.868.42"
$ws.Range("E12").Value = "  -9.33%  "
$ws.Range("D13").Value = "5.488"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "6.653"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "0.07202"
$ws.Range("E15").Value = "  +3.87%  "
$ws.Range("D16").Value = "82.38"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").Value = "1.034"
$ws.Range("E17").Value = "  +2.16%  "
$ws.Range("D18").Value = "0.000008990"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("D20").Value = "15.34"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").Value = "27.467.93"
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("D22").Value = "5.241"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").Value = "11.29"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").Value = "2.081.68"
$ws.Range("E24").Value = "  -9.45%  "
$ws.Range("D25").Value = "157.20"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("D26").Value = "1.920"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").Value = "18.62"
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("D28").Value = "5.239"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "1.920"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").Value = "116.55"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").Value = "0.09012"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").Value = "1.185"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").Value = "0.7546"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("D34").Value = "4.476"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").Value = "2.879"
$ws.Range("E35").Value = "  +2.48%  "
$ws.Range("D36").Value = "1.030"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("D38").Value = "0.01963"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("D39").Value = "0.05256"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.800"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5111"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").Value = "0.1659"
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").Value = "6.644"
$ws.Range("E43").Value = "  +2.17%  "
$ws.Range("D44").Value = "8.412"
$ws.Range("E44").Value = "  +2.00%  "
$ws.Range("D45").Value = "108.40"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("D46").Value = "10.49"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "1.700"
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.06395"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "0.4612"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("D50").Value = "1.846"
$ws.Range("E50").Value = "  +2.39%  "
$ws.Range("D51").Value = "39.08"
$ws.Range("E51").Value = "  +3.84%  "
